$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2436.3635
$ws.Range("J51").Value = 2428.5715
$ws.Range("L51").Value = 2428.5715
$ws.Range("N51").Value = -3396.5715

$ws.Range("H107").Value = 597.4
$ws.Range("I107").Value = 564.6875
$ws.Range("J107").Value = 728.25
$ws.Range("K107").Value = 564.6875
$ws.Range("L107").Value = 728.25
$ws.Range("M107").Value = 1355.3125
$ws.Range("N107").Value = -4568.25

$ws.Range("H112").Value = 6573.9287
$ws.Range("I112").Value = 35500
$ws.Range("J112").Value = 1752.9166
$ws.Range("K112").Value = 106500
$ws.Range("L112").Value = 5258.7498
$ws.Range("M112").Value = -105392
$ws.Range("N112").Value = -7474.7498

$ws.Range("H120").Value = 79761
$ws.Range("J120").Value = 79761
$ws.Range("L120").Value = 79761
$ws.Range("N120").Value = -89437

$ws.Range("H138").Value = 3456.6584
$ws.Range("I138").Value = 2676.6
$ws.Range("J138").Value = 3708.2903
$ws.Range("K138").Value = 8029.799999999999
$ws.Range("L138").Value = 11124.8709
$ws.Range("M138").Value = -2889.799999999999
$ws.Range("N138").Value = -21404.8709

$ws.Range("H141").Value = 2248.0833
$ws.Range("I141").Value = 1297.0385
$ws.Range("J141").Value = 4720.8
$ws.Range("K141").Value = 3891.1155
$ws.Range("L141").Value = 14162.4
$ws.Range("M141").Value = 1288.8845
$ws.Range("N141").Value = -24522.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 972.70966
$ws.Range("I2").Value = 990.8889
$ws.Range("J2").Value = 850
$ws.Range("K2").Value = 990.8889
$ws.Range("L2").Value = 850
$ws.Range("M2").Value = -877.8889
$ws.Range("N2").Value = -1076

$ws.Range("H32").Value = 23410.389
$ws.Range("I32").Value = 24757.377
$ws.Range("J32").Value = 8256.75
$ws.Range("K32").Value = 24757.377
$ws.Range("L32").Value = 8256.75
$ws.Range("M32").Value = -24470.377
$ws.Range("N32").Value = -8830.75

$ws.Range("H37").Value = 3017
$ws.Range("I37").Value = 3017
$ws.Range("K37").Value = 3017
$ws.Range("M37").Value = -2744

$ws.Range("H55").Value = 70000
$ws.Range("J55").Value = 70000
$ws.Range("L55").Value = 70000
$ws.Range("N55").Value = -70630

$ws.Range("H63").Value = 3998.7693
$ws.Range("I63").Value = 2914
$ws.Range("J63").Value = 4928.5713
$ws.Range("K63").Value = 2914
$ws.Range("L63").Value = 4928.5713
$ws.Range("M63").Value = -2228
$ws.Range("N63").Value = -6300.5713

$ws.Range("H66").Value = 3998.7693
$ws.Range("I66").Value = 2914
$ws.Range("J66").Value = 4928.5713
$ws.Range("K66").Value = 14570
$ws.Range("L66").Value = 24642.8565
$ws.Range("M66").Value = -11138
$ws.Range("N66").Value = -31506.8565

$ws.Range("H102").Value = 1280256.1
$ws.Range("I102").Value = 1686233
$ws.Range("J102").Value = 4328.5713
$ws.Range("K102").Value = 1686233
$ws.Range("L102").Value = 4328.5713
$ws.Range("M102").Value = -1684611
$ws.Range("N102").Value = -7572.5713

$ws.Range("H110").Value = 1635.1428
$ws.Range("I110").Value = 1684
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 1684
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 361
$ws.Range("N110").Value = -5090

$ws.Range("H116").Value = 972.70966
$ws.Range("I116").Value = 990.8889
$ws.Range("J116").Value = 850
$ws.Range("K116").Value = 990.8889
$ws.Range("L116").Value = 850
$ws.Range("M116").Value = 1303.1111
$ws.Range("N116").Value = -5438

$ws.Range("H122").Value = 4809445
$ws.Range("I122").Value = 1826.0952
$ws.Range("K122").Value = 5478.2856
$ws.Range("M122").Value = -3028.2856

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 972.70966
$ws.Range("I3").Value = 990.8889
$ws.Range("J3").Value = 850
$ws.Range("K3").Value = 990.8889
$ws.Range("L3").Value = 850
$ws.Range("M3").Value = -876.8889
$ws.Range("N3").Value = -1078

$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").ClearContents()
$ws.Range("N35").Value = 0

$ws.Range("H82").Value = 21261.8
$ws.Range("I82").Value = 4618.375
$ws.Range("J82").Value = 40282.855
$ws.Range("K82").Value = 4618.375
$ws.Range("L82").Value = 40282.855
$ws.Range("M82").Value = -4235.375
$ws.Range("N82").Value = -41048.855

$ws.Range("H85").Value = 21261.8
$ws.Range("I85").Value = 4618.375
$ws.Range("J85").Value = 40282.855
$ws.Range("K85").Value = 4618.375
$ws.Range("L85").Value = 40282.855
$ws.Range("M85").Value = -3292.375
$ws.Range("N85").Value = -42934.855

$ws.Range("H105").Value = 1206136
$ws.Range("I105").Value = 1648881
$ws.Range("J105").Value = 4399.5713
$ws.Range("K105").Value = 1648881
$ws.Range("L105").Value = 4399.5713
$ws.Range("M105").Value = -1647134
$ws.Range("N105").Value = -7893.5713

$ws.Range("H132").Value = 50589.5
$ws.Range("J132").Value = 50589.5
$ws.Range("L132").Value = 50589.5
$ws.Range("N132").Value = -60709.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1961.4445
$ws.Range("I16").Value = 1164.7778
$ws.Range("J16").Value = 2758.111
$ws.Range("K16").Value = 1164.7778
$ws.Range("L16").Value = 2758.111
$ws.Range("M16").Value = -877.7778000000001
$ws.Range("N16").Value = -3332.111

$ws.Range("H31").Value = 564195.0600000001
$ws.Range("I31").Value = 4356.5
$ws.Range("J31").Value = 976707.7
$ws.Range("K31").Value = 4356.5
$ws.Range("L31").Value = 976707.7
$ws.Range("M31").Value = -4061.5
$ws.Range("N31").Value = -977297.7

$ws.Range("H34").Value = 564195.0600000001
$ws.Range("I34").Value = 4356.5
$ws.Range("J34").Value = 976707.7
$ws.Range("K34").Value = 4356.5
$ws.Range("L34").Value = 976707.7
$ws.Range("M34").Value = -4154.5
$ws.Range("N34").Value = -977111.7

$ws.Range("H94").Value = 253.15384
$ws.Range("J94").Value = 301.2
$ws.Range("L94").Value = 301.2
$ws.Range("N94").Value = -1203.2

$ws.Range("H113").Value = 1961.4445
$ws.Range("I113").Value = 1164.7778
$ws.Range("J113").Value = 2758.111
$ws.Range("K113").Value = 1164.7778
$ws.Range("L113").Value = 2758.111
$ws.Range("M113").Value = 1005.2222
$ws.Range("N113").Value = -7098.111

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2600.5
$ws.Range("I80").Value = 2501
$ws.Range("J80").Value = 2700
$ws.Range("K80").Value = 7503
$ws.Range("L80").Value = 8100
$ws.Range("M80").Value = -6567
$ws.Range("N80").Value = -9972

$ws.Range("H83").Value = 2600.5
$ws.Range("I83").Value = 2501
$ws.Range("J83").Value = 2700
$ws.Range("K83").Value = 22509
$ws.Range("L83").Value = 24300
$ws.Range("M83").Value = -17829
$ws.Range("N83").Value = -33660

$ws.Range("H107").Value = 565.7564
$ws.Range("I107").Value = 324.08334
$ws.Range("J107").Value = 1371.3334
$ws.Range("K107").Value = 972.2500200000001
$ws.Range("L107").Value = 4114.0002
$ws.Range("M107").Value = 947.7499799999999
$ws.Range("N107").Value = -7954.0002

$ws.Range("H122").Value = 1303.9
$ws.Range("I122").Value = 550
$ws.Range("J122").Value = 1492.375
$ws.Range("K122").Value = 4950
$ws.Range("L122").Value = 13431.375
$ws.Range("M122").Value = -2500
$ws.Range("N122").Value = -18331.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2105.0527
$ws.Range("I126").Value = 1510.8334
$ws.Range("J126").Value = 3123.7144
$ws.Range("K126").Value = 4532.5002
$ws.Range("L126").Value = 9371.143199999999
$ws.Range("M126").Value = -2062.5002
$ws.Range("N126").Value = -14311.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3435.0454
$ws.Range("I40").Value = 3307.6924
$ws.Range("K40").Value = 3307.6924
$ws.Range("M40").Value = -3171.6924

$ws.Range("H46").Value = 984.5
$ws.Range("I46").Value = 720.1111
$ws.Range("K46").Value = 720.1111
$ws.Range("M46").Value = -532.1111

$ws.Range("H61").Value = 23667.549
$ws.Range("I61").Value = 24675.72
$ws.Range("J61").Value = 19466.834
$ws.Range("K61").Value = 24675.72
$ws.Range("L61").Value = 19466.834
$ws.Range("M61").Value = -24473.72
$ws.Range("N61").Value = -19870.834

$ws.Range("H113").Value = 23667.549
$ws.Range("I113").Value = 24675.72
$ws.Range("J113").Value = 19466.834
$ws.Range("K113").Value = 24675.72
$ws.Range("L113").Value = 19466.834
$ws.Range("M113").Value = -22505.72
$ws.Range("N113").Value = -23806.834

$ws.Range("H122").Value = 7296.4443
$ws.Range("I122").Value = 6999.9546
$ws.Range("J122").Value = 8601
$ws.Range("K122").Value = 20999.8638
$ws.Range("L122").Value = 25803
$ws.Range("M122").Value = -18549.8638
$ws.Range("N122").Value = -30703

$ws.Range("H136").Value = 2661.5088
$ws.Range("I136").Value = 1998.138
$ws.Range("J136").Value = 3348.5715
$ws.Range("K136").Value = 5994.414
$ws.Range("L136").Value = 10045.7145
$ws.Range("M136").Value = -3444.414
$ws.Range("N136").Value = -15145.7145

$ws.Range("H139").Value = 51666.668
$ws.Range("J139").Value = 51666.668
$ws.Range("L139").Value = 51666.668
$ws.Range("N139").Value = -61946.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 772.5
$ws.Range("I113").Value = 400.375
$ws.Range("J113").Value = 1070.2
$ws.Range("K113").Value = 1201.125
$ws.Range("L113").Value = 3210.6
$ws.Range("M113").Value = 968.875
$ws.Range("N113").Value = -7550.6
